$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.397.09'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.938.36'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7681'
$ws.Range("E5").Value = '  +8.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '247.99'
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.99'
$ws.Range("E8").Value = '  +2.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3203'
$ws.Range("E9").Value = '  -2.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07114'
$ws.Range("E10").Value = '  -2.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7845'
$ws.Range("E11").Value = '  -2.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08017'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.938.31'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.382'
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.18'
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.55'
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.402.12'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '258.22'
$ws.Range("E18").Value = '  +1.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008019'
$ws.Range("E19").Value = '  -2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.855'
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.191.59'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9988'
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.761'
$ws.Range("E24").Value = '  -3.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.617'
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.56'
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.17'
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1332'
$ws.Range("E28").Value = '  +3.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.295'
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("E30").Value = '  +1.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.526'
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.442'
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.158'
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05207'
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.282'
$ws.Range("E35").Value = '  +1.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7503'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.778'
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01972'
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.807'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.57'
$ws.Range("E40").Value = '  -0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.466'
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4532'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.973'
$ws.Range("E43").Value = '  -2.08%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8343'
$ws.Range("E45").Value = '  -1.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.44'
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.825'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.539'
$ws.Range("E48").Value = '  +1.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '986.70'
$ws.Range("E49").Value = '  +11.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.42'
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4169'
$ws.Range("E51").Value = '  -0.25%  '
